# 06组项目计划表.xlsx - "Add files via upload" re-edit
#
# The underlying change (once the shared-strings / cellXfs index churn that
# Excel/WPS produces on every re-save is factored out) is a simple content
# fix on Sheet1: the two "第八周" ("week 8") date headers were corrected to
# read "第七周" ("week 7"), and the active selection moved from B34 to I27.
#
# A19 ("日期：2018.10.10 第七周周三") is textually identical before/after, so
# no write is required there.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 29 header: "日期：2018.10.15 第八周周一" -> "日期：2018.10.15 第七周周一"
$ws.Range("A29").Value = "日期：2018.10.15 第七周周一"

# Row 39 header: "日期：2018.10.15 第八周周三" -> "日期：2018.10.15 第七周周三"
$ws.Range("A39").Value = "日期：2018.10.15 第七周周三"

# Selection moved to I27
$null = $ws.Range("I27").Select()
